# SMLR Yearly Financials update: add FY2018 (period ending 2018-12-31) column
# by inserting a new column before column D and populating it with the
# newest-year figures across the Income Statement, Balance Sheet and
# Cash Flow Statement tables. Everything that used to live in D:K shifts
# right to E:L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (FY2018) values keyed by row number - covers every row that carries
# data in columns D:K in all three statements.
$newDValues = @{
    7 = 43465
    8 = 21500
    9 = 2700
    10 = 18800
    12 = 2100
    13 = 0
    14 = 0
    15 = 0
    17 = 16100
    18 = 5300
    20 = 0
    21 = 5800
    22 = 300
    23 = 5000
    24 = 0
    25 = 0
    26 = 5000
    27 = 5000
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 5000
    34 = 0
    35 = 5000
    38 = 43465
    41 = 3300
    42 = 0
    43 = 2800
    44 = 0
    45 = 200
    46 = 6200
    47 = 0
    48 = 1500
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 7700
    57 = 300
    58 = 0
    59 = 3200
    60 = 3500
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 3500
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -21400
    73 = 0
    74 = 0
    75 = 0
    76 = 4200
    77 = 0
    80 = 43465
    81 = 5000
    83 = 500
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 4700
    91 = -800
    92 = 0
    93 = 0
    94 = -800
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = -2000
    101 = 0
    102 = 1800
}

# Rows that end with an "NA" marker (shared text "NA") rather than a number
# in what is now column K (was column J before the shift).
$naRows = @(14, 24, 62)

# 1. Insert a blank column before D; this shifts old D:K -> E:L.
$ws.Columns("D:D").Insert()

# 2. The new column D defaults to the plain/general column style, so pull
#    the real number formats (date format for the header rows, #,##0 for
#    the data rows) back from column E, which now holds the formatting
#    that used to belong to column D.
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3. Populate the new column D with the FY2018 figures.
foreach ($row in $newDValues.Keys) {
    $ws.Cells.Item($row, 4).Value2 = $newDValues[$row]
}

# 4. Restore the "NA" shared-text marker in the new column K for the rows
#    where the last historical year had no data (this mirrors what
#    shifted in from the old column J).
foreach ($row in $naRows) {
    $ws.Cells.Item($row, 11).Value2 = "NA"
}

# 5. Dimension / used range now reaches column L instead of K.
$ws.Range("A5:L102").Select() | Out-Null

$wb.Save()
